$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '332.26'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '1.93%'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '45.95'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '3.73%'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.704'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '3.67%'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.08358'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '4.51%'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '2.052'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '1.35%'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.9743'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '2.98%'
$ws.Range("B8").Value = 'BTSEToken'
$ws.Range("C8").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '2.592'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '1.06%'
$ws.Range("B9").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C9").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.1153'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '1.40%'
$ws.Range("B10").Value = 'WazirX'
$ws.Range("C10").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.1942'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '5.51%'
$ws.Range("B11").Value = 'MCDex'
$ws.Range("C11").Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '10.43'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '-14.18%'
$ws.Range("B12").Value = 'MandalaExchangeToken'
$ws.Range("C12").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.1005'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '3.19%'
$ws.Range("B13").Value = 'BitrueCoin'
$ws.Range("C13").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.04606'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '1.03%'
$ws.Range("B14").Value = 'BitMartToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.1059'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '-0.47%'
$ws.Range("B15").Value = 'BitForexToken'
$ws.Range("C15").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.001293'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '1.79%'
$ws.Range("B16").Value = 'TigerCash'
$ws.Range("C16").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.006097'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '3.51%'
$ws.Range("B17").Value = 'LEO'
$ws.Range("C17").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.368'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '0.11%'
$ws.Range("B18").Value = 'GateToken'
$ws.Range("C18").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '4.464'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '3.78%'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.3351'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '-3.68%'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.1390'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '-1.37%'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.2593'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '1.87%'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.04183'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '2.70%'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.001307'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '5.32%'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.004715'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '9.78%'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.0001281'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '7.75%'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.0003741'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '-0.05%'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02772'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '9.44%'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.05842'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '5.51%'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.007743'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '2.86%'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.1438'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '3.54%'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.007194'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '-5.26%'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.001975'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '-1.90%'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.008175'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '-2.49%'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.00007191'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '1.28%'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '0.16%'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0005802'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '-0.15%'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.003488'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '-1.21%'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.003499'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '52.02%'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.00002101'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '0.16%'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0002001'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '0.16%'

Write-Output "Applied 99 cell changes"